# Daily attendance processing - 2026-01-18 17:33:26
# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in the "Recorded By" column (column G) of the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
